$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions (C1:F1) ---
$ws.Range("C1").Value = "uploadConfigPath"
$ws.Range("D1").Value = "productNum"
$ws.Range("E1").Value = "productDesc"
$ws.Range("F1").Value = "productOpt"

# --- Row 2 additions (C2, D2, E2) ---
$ws.Range("C2").Value = "DESSTEPS_335038_CPQ Encore_US9400_04_config.xls"

# Columns D and E (rows 2-12) are formatted as Text ("@") before entry,
# matching the workbook's applied "Text" cell style. Row 11 in column E is
# the exception -- it keeps the default General format (see below).
$ws.Range("D2:D12").NumberFormat = "@"
$ws.Range("E2:E10").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"

$ws.Range("D2").Value = "727021-B21"
$ws.Range("E2").Value = "HP BL460c Gen9 10Gb/20Gb FLB CTO Blade"

$ws.Range("D3").Value = "726991-L21"
$ws.Range("E3").Value = "HP BL460c Gen9 E5-2650v3 FIO Kit"

$ws.Range("D4").Value = "726991-B21"
$ws.Range("E4").Value = "HP BL460c Gen9 E5-2650v3 Kit"

$ws.Range("D5").Value = "726991-B21"
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = "0D1"

$ws.Range("D6").Value = "726722-B21"
$ws.Range("E6").Value = "HP 32GB 4Rx4 PC4-2133P-L Kit"

$ws.Range("D7").Value = "726722-B21"
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = "0D1"

$ws.Range("D8").Value = "700764-B21"
$ws.Range("E8").Value = "HP FlexFabric 20Gb 2P 650FLB FIO Adptr"

$ws.Range("D9").Value = "H1K92A3"
$ws.Range("E9").Value = "HPE 3Y Proactive Care 24x7 SVC"

$ws.Range("D10").Value = "H1K92A3"
$ws.Range("E10").Value = "HP BL460c Gen9 Server Blade HW Supp"
$ws.Range("F10").Value = "TT8"

$ws.Range("D11").Value = "HA114A1"
$ws.Range("E11").Value = "HP CP Installation & Startup"
# E11 keeps the default (General) number format but uses the dark-grey font
# color used elsewhere in the source workbook.
$ws.Range("E11").Font.Color = 2236962

$ws.Range("D12").Value = "HA114A1"
$ws.Range("E12").Value = "HP C Class Server Blade Startup SVC"
$ws.Range("F12").Value = "5CY"

# --- Column widths for the new columns (C=84, D=17.86, E=38.71 "bestFit", F=26.71 chars) ---
$ws.Range("C1").ColumnWidth = 83.16666666666667
$ws.Range("D1").ColumnWidth = 17
$ws.Range("E1").ColumnWidth = 37.833333333333336
$ws.Range("F1").ColumnWidth = 25.833333333333332

# --- View / selection / print setup ---
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollColumn = 4
$ws.PageSetup.Orientation = 1
